# Add a new MEDIA_TYPE column (E) to the PROMPTS_CONFIG sheet, populate it
# with "TEXT" for every existing data row (including two previously-blank
# rows that become real rows once column E has content), resize columns D
# and E, and restore the view state (selections) on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- PROMPTS_CONFIG (sheet1): add column E -------------------------------
$ws1.Range("E1").Value = "MEDIA_TYPE"

for ($r = 2; $r -le 11; $r++) {
    $ws1.Cells.Item($r, 5).Value = "TEXT"
}

# Column D loses its shared bestFit with C and gets an explicit width;
# column E gets its own explicit width too.
$ws1.Range("D1").ColumnWidth = 37.166666666666664
$ws1.Range("E1").ColumnWidth = 10.654296875

# --- RICH_CARDS_CONFIG (sheet2): update remembered selection -------------
# Selecting on sheet2 activates it, so do this before re-selecting sheet1
# so that PROMPTS_CONFIG ends up as the active tab again.
$null = $ws2.Range("B2").Select()

# --- restore PROMPTS_CONFIG as the active sheet/selection -----------------
$null = $ws1.Range("C14").Select()
